$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.266.84"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "2.633.59"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.94%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.572"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.94%  "
$ws.Range("D9").Value = "2.641.36"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("E10").Value = "  -5.17%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.341"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").Value = "3.090.95"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").Value = "60.247.81"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").Value = "2.668.63"
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.419"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.60%  "
$ws.Range("E28").Value = "  -4.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("E35").Value = "  -5.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.916"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.69%  "
$ws.Range("E37").Value = "  -5.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.863"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "293.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.85%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.630"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.74%  "
$ws.Range("D51").Value = "1.961.82"
$ws.Range("E51").Value = "  -0.09%  "
